$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "299.57"
Set-TextValue $ws.Range("E2") "1.81%"
Set-TextValue $ws.Range("D3") "31.23"
Set-TextValue $ws.Range("E3") "-0.10%"
Set-TextValue $ws.Range("D4") "5.131"
Set-TextValue $ws.Range("E4") "0.80%"
Set-TextValue $ws.Range("D5") "0.08118"
Set-TextValue $ws.Range("E5") "10.21%"
Set-TextValue $ws.Range("D6") "2.498"
Set-TextValue $ws.Range("E6") "53.14%"
Set-TextValue $ws.Range("D7") "7.852"
Set-TextValue $ws.Range("E7") "2.19%"
Set-TextValue $ws.Range("D8") "3.842"
Set-TextValue $ws.Range("E8") "2.29%"
Set-TextValue $ws.Range("D9") "0.9096"
Set-TextValue $ws.Range("E9") "-1.58%"
Set-TextValue $ws.Range("D10") "0.1712"
Set-TextValue $ws.Range("E10") "2.60%"
Set-TextValue $ws.Range("D11") "0.07277"
Set-TextValue $ws.Range("E11") "1.29%"
Set-TextValue $ws.Range("D12") "0.07983"
Set-TextValue $ws.Range("E12") "0.49%"
Set-TextValue $ws.Range("E13") "1.05%"
Set-TextValue $ws.Range("D14") "0.09967"
Set-TextValue $ws.Range("E14") "0.75%"
Set-TextValue $ws.Range("D15") "0.001502"
Set-TextValue $ws.Range("E15") "0.76%"
Set-TextValue $ws.Range("D16") "0.005993"
Set-TextValue $ws.Range("E16") "-3.49%"
Set-TextValue $ws.Range("E17") "1.16%"
Set-TextValue $ws.Range("D19") "0.3256"
Set-TextValue $ws.Range("E19") "-0.69%"
Set-TextValue $ws.Range("D20") "0.1347"
Set-TextValue $ws.Range("E20") "0.87%"
Set-TextValue $ws.Range("D21") "4.593"
Set-TextValue $ws.Range("E21") "0.54%"
Set-TextValue $ws.Range("D22") "0.1603"
Set-TextValue $ws.Range("D23") "0.04585"
Set-TextValue $ws.Range("E23") "-0.95%"
Set-TextValue $ws.Range("D24") "0.001265"
Set-TextValue $ws.Range("E24") "3.84%"
Set-TextValue $ws.Range("D25") "0.004441"
Set-TextValue $ws.Range("E25") "0.48%"
Set-TextValue $ws.Range("E26") "-9.15%"
Set-TextValue $ws.Range("D27") "0.0003436"
Set-TextValue $ws.Range("E27") "83.06%"
Set-TextValue $ws.Range("D39") "0.01815"
Set-TextValue $ws.Range("E39") "7.72%"
Set-TextValue $ws.Range("D40") "0.04545"
Set-TextValue $ws.Range("E40") "3.29%"
Set-TextValue $ws.Range("D41") "0.007165"
Set-TextValue $ws.Range("E41") "0.58%"
Set-TextValue $ws.Range("E42") "1.26%"
Set-TextValue $ws.Range("D43") "0.002244"
Set-TextValue $ws.Range("E43") "6.77%"
Set-TextValue $ws.Range("D44") "0.01052"
Set-TextValue $ws.Range("E44") "-4.56%"
Set-TextValue $ws.Range("D45") "0.00006282"
Set-TextValue $ws.Range("E45") "4.61%"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("E46") "0.18%"
Set-TextValue $ws.Range("B47") "CoinbaseStockToken"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws.Range("D47") "0.006412"
Set-TextValue $ws.Range("E47") "-41.76%"
Set-TextValue $ws.Range("B48") "BOLO"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws.Range("D48") "0.8206"
Set-TextValue $ws.Range("E48") "-57.22%"
Set-TextValue $ws.Range("D49") "0.00002104"
Set-TextValue $ws.Range("E49") "0.18%"
Set-TextValue $ws.Range("D50") "0.0002004"
Set-TextValue $ws.Range("E50") "0.25%"
